$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "72.770.64"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +4.98%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "4.056.75"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +4.07%  "

$ws.Range("E4").Value = "  +0.17%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "518.55"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.93%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.99"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.64%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.738"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +20.52%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "4.047.66"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.27%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.14%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.772"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +7.49%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.176"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.68%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000327"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.88%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "47.74"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +13.55%  "

$ws.Range("E14").Value = "  +8.92%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.707.56"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.15%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.060.53"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.27%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "21.19"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +7.10%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.08"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.65%  "

$ws.Range("E19").Value = "  -0.39%  "

$ws.Range("E20").Value = "  -1.17%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.607.73"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.86%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "443.62"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.19%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "104.91"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +19.01%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.56"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.43%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "14.77"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.69%  "

$ws.Range("E26").Value = "  -0.83%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.44"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.77%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.01"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.92%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "37.83"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.98%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.82"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.66%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.28"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +15.94%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "13.65"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.00%  "

$ws.Range("E33").Value = "  +3.35%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "681.53"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.43%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.81"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +14.57%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "67.08"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.45%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "42.63"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +6.81%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0861"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.36%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.428"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.63%  "

$ws.Range("E40").Value = "  +7.41%  "

$ws.Range("E41").Value = "  +1.70%  "

$ws.Range("E42").Value = "  -0.02%  "

$ws.Range("E43").Value = "  +3.57%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.999"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.17%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.28"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.47%  "

$ws.Range("E46").Value = "  +12.86%  "

$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.50"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.69%  "

$ws.Range("B48").Value = "Fetch.AI"
$ws.Range("C48").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.69"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.03%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.07"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.19%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.18"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +7.59%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.32"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.25%  "
